$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy number formats from the last existing row (24) down into the new
# row (25) first - a plain values+formats paste would make Excel derive
# an ad-hoc date numFmt instead of reusing the workbook's existing style.
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)

# Copy over the values that stay identical to the row above (Date,
# Hand in, By) with a values-only paste, so e.g. "TRUE" round-trips as
# the existing shared string instead of becoming a COM boolean.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4163)
$ws.Range("D24").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E24").Copy()
$ws.Range("E25").PasteSpecial(-4163)

# New submission details for this row.
$ws.Range("B25").Value = "221126_cat_test_2"
$ws.Range("C25").Value = "pineapple_pizza"

# Extend the table (ListObject) to include the new row.
$table = $ws.ListObjects.Item("Tabelle2")
$table.Resize($ws.Range("A1:E25"))

$ws.Range("E26").Select()
